$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library")

$ws.Range("B2").Value = "nucleic acid library construction protocol"
$ws.Range("C2").Value = "EFO"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/EFO_0004184"
$ws.Range("E2").Value = "library_construction.txt"
